$d = $word.ActiveDocument
$p19 = $d.Paragraphs.Item(19)
$p25 = $d.Paragraphs.Item(25)
$startPos = $p19.Range.Start
$endPos = $p25.Range.End
$rng = $d.Range($startPos, $endPos)

$xmlContent = @"
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">{label: &#8220;node1&#8221;, children: </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">[{ </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>edgeLabel</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">: &#8220;to_node_2&#8221;, label: node_2 }, { </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>edgeLabel</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">: &#8220;to_node_3&#8221;, label: &#8220;node_3&#8221;, children [</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>&#8230;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>]  }]}</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Each node is exported as a label attribute indicating the attribute of the dataset we are considering at that node and, if it is not the root, an </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>edgeLabel</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> attribute </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">indicating which branch we have to follow to reach the current node from its parent, therefore the answer to the parent&#8217;s question. </w:t>
  </w:r>
  <w:r>
    <w:t>Moreover, a children attribute, which is an array of nodes, could be present</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> for non-leaf</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve"> nodes</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">The </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">JSON input file is the only parameter required to the bot interpreter application which then creates a RESTful server on the local machine and a Microsoft </w:t>
  </w:r>
  <w:r>
    <w:t>&#8220;</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>botbuilder</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>&#8221;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> object. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">The application parses then the JSON specification file and presents the user a welcoming message. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">When the user is ready to start, it creates a dialog with the root node as parameter. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>The general strategy of the interpreter is to create a new dialog for each node, replacing the older one</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> to keep a low memory usage</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">The dialog reads the attribute that needs to be asked to the user and first checks if an answer to the current question was already given previously. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">This is important especially for numeric </w:t>
  </w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">attributes. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">C4.5 creates for each numeric attribute a binary split, where a branch contains all instances where the attribute is less or equal to the split point, while the other branch all ones with a greater value. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">In some cases, value ranges are needed to precisely decide between different cases and </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">so the algorithm asks more than one time the same &#8220;question&#8221; with different split points. This could happen at different levels of the decision tree but also on consecutive levels. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">If the bot interpreter application </w:t>
  </w:r>
  <w:r>
    <w:t>would</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> not memorize the questions and answers </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">pairs, it will </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">ask more than one time the same question, forcing the user to answer multiple times. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">For this reason, every time a new question is posed to the user, it is memorized along with the collected answer in a map. </w:t>
  </w:r>
  <w:r>
    <w:t>We are so able to check if an answer has already been asked</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>and, in that case,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">retrieve the answer and use it to decide which branch to follow. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">On the other hand, if the question has not already been asked, we pose it to the user. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">If the attribute </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">is a numerical attribute, the user is just prompted with the question and is able </w:t>
  </w:r>
  <w:r>
    <w:t>to answer with whatever is a number</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> (bot integers and floating points). </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Differently, if the attribute is categorical and therefore it can </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">only take some specific values, </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">the application prompts the questions and a list of clickable buttons representing the accepted values. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">The user can so directly click on one of the given options or reply with a message containing one of the options&#8217; text. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">When the interpreter reaches a leaf node, so a node with the children attribute missing, it presents the </w:t>
  </w:r>
  <w:r>
    <w:t>answer to the user and exits</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> dumping on the server console all questions and answers retrieved that were memorized in the map. </w:t>
  </w:r>
</w:p>

</pkg:xmlData>
"@

$rng.InsertXML($xmlContent)
Write-Host "Replacement done"
